$wb = $excel.ActiveWorkbook
$wsProperty = $wb.Worksheets.Item("Property")
$wsPosList  = $wb.Worksheets.Item("Record_PosList")

# ---------------------------------------------------------------------------
# 1) Swap the header row (row 1) and the sample/type row (row 2) for columns
#    L:O on the Record_PosList sheet. Before the edit, row 1 held the
#    (incorrect) "float" placeholders and row 2 held the real X/Y/Z/StayTime
#    headers; after the edit row 1 holds the headers and row 2 holds the
#    "float" type markers.
# ---------------------------------------------------------------------------
$cols = @("L", "M", "N", "O")
$row1Values = @{}
$row2Values = @{}
foreach ($col in $cols) {
    $row1Values[$col] = $wsPosList.Range("${col}1").Value2
    $row2Values[$col] = $wsPosList.Range("${col}2").Value2
}
foreach ($col in $cols) {
    $wsPosList.Range("${col}1").Value = $row2Values[$col]
    $wsPosList.Range("${col}2").Value = $row1Values[$col]
}

# ---------------------------------------------------------------------------
# 2) Move the column comments from row 1 down to row 2, updating the text of
#    some of them along the way (matches the target commit exactly).
# ---------------------------------------------------------------------------
$wsPosList.Range("L1").Comment.Delete()
$wsPosList.Range("M1").Comment.Delete()
$wsPosList.Range("N1").Comment.Delete()
$wsPosList.Range("O1").Comment.Delete()

$wsPosList.Range("L2").AddComment("强化等级")
$wsPosList.Range("M2").AddComment("强化等级")
$wsPosList.Range("N2").AddComment("强化等级")
$wsPosList.Range("O2").AddComment("镶嵌宝石，逗号分隔")

# ---------------------------------------------------------------------------
# 3) Update the view state: Record_PosList becomes the active tab (and its
#    prior "tabSelected" on Property is cleared), with new selections on
#    both sheets.
# ---------------------------------------------------------------------------
[void]$wsProperty.Activate()
[void]$wsProperty.Range("J33").Select()

[void]$wsPosList.Activate()
[void]$wsPosList.Range("O8").Select()
